$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F48").Value = 31
$ws.Range("G48").Value = 6099.87
$ws.Range("F62").Value = 3
$ws.Range("G62").Value = 168.3
$ws.Range("B85").Value = 157259.6
$ws.Range("F98").Value = 8
$ws.Range("G98").Value = 3017.52
$ws.Range("B103").Value = 26437.54
$ws.Range("F143").Value = 40
$ws.Range("G143").Value = 1594.4
$ws.Range("B159").Value = 72110.71000000001
$ws.Range("F169").Value = 65
$ws.Range("G169").Value = 3216.2
$ws.Range("F171").Value = 117
$ws.Range("G171").Value = 4951.44
$ws.Range("B180").Value = 36585.35
$ws.Range("F185").Value = 31
$ws.Range("G185").Value = 2024.3
$ws.Range("B198").Value = 44340.52
$ws.Range("F218").Value = 24
$ws.Range("G218").Value = 1878.24
$ws.Range("F219").Value = 206
$ws.Range("G219").Value = 6392.18
$ws.Range("B228").Value = 35293.68
$ws.Range("F287").Value = 18
$ws.Range("G287").Value = 4153.68
$ws.Range("B298").Value = 109713.29
$ws.Range("F307").Value = 3
$ws.Range("G307").Value = 451.83
$ws.Range("B310").Value = 30046.6
$ws.Range("F322").Value = 23
$ws.Range("G322").Value = 1063.75
$ws.Range("F329").Value = 201
$ws.Range("G329").Value = 6693.3
$ws.Range("F343").Value = 18
$ws.Range("G343").Value = 843.66
$ws.Range("F346").Value = 2
$ws.Range("G346").Value = 164.32
$ws.Range("B349").Value = 141284.25
$ws.Range("F382").Value = 27
$ws.Range("G382").Value = 2965.14
$ws.Range("F389").Value = 45
$ws.Range("G389").Value = 5708.7
$ws.Range("F399").Value = 122
$ws.Range("G399").Value = 7161.4
$ws.Range("F402").Value = 10
$ws.Range("G402").Value = 544.8
$ws.Range("F419").Value = 283
$ws.Range("G419").Value = 11653.94
$ws.Range("B423").Value = 147581.29
$ws.Range("F425").Value = 16
$ws.Range("G425").Value = 2937.28
$ws.Range("F436").Value = 3
$ws.Range("G436").Value = 670.8
$ws.Range("B437").Value = 21066.66
$ws.Range("F440").Value = 0
$ws.Range("G440").Value = 0
$ws.Range("B441").Value = 0
$ws.Range("F469").Value = 31
$ws.Range("G469").Value = 2882.69
$ws.Range("B481").Value = 43776.47
$ws.Range("F495").Value = 29
$ws.Range("G495").Value = 4813.71
$ws.Range("B497").Value = 38789.84
$ws.Range("F509").Value = 23
$ws.Range("G509").Value = 1536.17
$ws.Range("B511").Value = 36730.58
$ws.Range("F516").Value = 179
$ws.Range("G516").Value = 12288.35
$ws.Range("F526").Value = 771
$ws.Range("G526").Value = 74478.60000000001
$ws.Range("B532").Value = 153238.49
$ws.Range("F560").Value = 109
$ws.Range("G560").Value = 5509.95
$ws.Range("B567").Value = 51388.05
$ws.Range("F576").Value = 0
$ws.Range("G576").Value = 0
$ws.Range("B583").Value = 38987.74
$ws.Range("F597").Value = 122
$ws.Range("G597").Value = 7521.3
$ws.Range("F606").Value = 166
$ws.Range("G606").Value = 5778.46
$ws.Range("B610").Value = 56530.83
$ws.Range("F621").Value = 240
$ws.Range("G621").Value = 14572.8
$ws.Range("B623").Value = 55667
$ws.Range("C623").Value = 'NES-Maggi Atta Noodles Masala 290G'
$ws.Range("D623").Value = 85.76000000000001
$ws.Range("E623").Value = 97.25
$ws.Range("F623").Value = 71
$ws.Range("G623").Value = 6088.96
$ws.Range("B624").Value = 49151
$ws.Range("C624").Value = 'NES-MAGGI Atta Noodles Masala 290g'
$ws.Range("D624").Value = 78.09999999999999
$ws.Range("E624").Value = 88.58
$ws.Range("F624").Value = 1
$ws.Range("G624").Value = 78.09999999999999
$ws.Range("B638").Value = 147362.43
$ws.Range("F664").Value = 18
$ws.Range("G664").Value = 1360.08
$ws.Range("F666").Value = 39
$ws.Range("G666").Value = 2701.53
$ws.Range("B667").Value = 26081.14
$ws.Range("F672").Value = 69
$ws.Range("G672").Value = 18360.21
$ws.Range("F673").Value = 143
$ws.Range("G673").Value = 3729.44
$ws.Range("F678").Value = 16
$ws.Range("G678").Value = 751.36
$ws.Range("B688").Value = 87143.14999999999
$ws.Range("F715").Value = 46
$ws.Range("G715").Value = 6005.3
$ws.Range("F717").Value = 44
$ws.Range("G717").Value = 1196.8
$ws.Range("F718").Value = 128
$ws.Range("G718").Value = 3481.6
$ws.Range("F719").Value = 106
$ws.Range("G719").Value = 2883.2
$ws.Range("B720").Value = 30633.24
$ws.Range("F825").Value = 21
$ws.Range("G825").Value = 782.88
$ws.Range("B837").Value = 194257.32
$ws.Range("F843").Value = 68
$ws.Range("G843").Value = 7399.08
$ws.Range("F862").Value = 14
$ws.Range("G862").Value = 660.9400000000001
$ws.Range("F865").Value = 96
$ws.Range("G865").Value = 4789.44
$ws.Range("B867").Value = 203127.27
$ws.Range("B923").Value = 2601490.93
$ws.Range("B924").Value = 2601490.93
